$d = $word.ActiveDocument

# "Updated Batch FTA template to judge magistrate."
# Replace the judge signature line:
#   "Judge Marianne T. Hemmeter / Judge Kyle E. Rohrer"
# with:
#   "Judge / Magistrate"

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("Judge Marianne T. Hemmeter", $true, $false, $false, $false, $false, `
               $true, 1, $false, "Judge ", 2)

$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Replacement.ClearFormatting()
$find2.Execute(" / Judge Kyle E. Rohrer", $true, $false, $false, $false, $false, `
               $true, 1, $false, "/ Magistrate", 2)
